# Insert a new row for Malawi (MWI) before the existing Malaysia row (row 76),
# shifting Malaysia and all subsequent countries down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = 76

# Insert a new blank row at position 76; existing row 76 (Malaysia) and below shift to 77+
$ws.Rows.Item($targetRow).Insert()

# Populate the new row with Malawi's data, following the same column layout as every
# other data row in the sheet:
# A=freqCode, B=date, C=refYear, D=refMonth, E=reporterCode, F=reporterISO,
# G=reporterDesc, H=flowDesc, I=partnerCode, J=partnerISO, K=partnerDesc, L=primaryValue
$ws.Cells.Item($targetRow, 1).Value2 = "A"
$ws.Cells.Item($targetRow, 2).Value2 = 44197
$ws.Cells.Item($targetRow, 2).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item($targetRow, 3).Value2 = 2021
$ws.Cells.Item($targetRow, 4).Value2 = 52
$ws.Cells.Item($targetRow, 5).Value2 = 454
$ws.Cells.Item($targetRow, 6).Value2 = "MWI"
$ws.Cells.Item($targetRow, 7).Value2 = "Malawi"
$ws.Cells.Item($targetRow, 8).Value2 = "Export"
$ws.Cells.Item($targetRow, 9).Value2 = 0
$ws.Cells.Item($targetRow, 10).Value2 = "W00"
$ws.Cells.Item($targetRow, 11).Value2 = "World"
$ws.Cells.Item($targetRow, 12).Value2 = 1009460778.534
